# Auto-generated update of market/profit figures per the commit diff.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H-N) for the
# affected Leve rows across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 5774.5
$ws.Cells.Item(32, 10).Value = 5774.5
$ws.Cells.Item(32, 12).Value = 5774.5
$ws.Cells.Item(32, 14).Value = -6426.5

$ws.Cells.Item(70, 8).Value = 2300.889
$ws.Cells.Item(70, 9).Value = 1714.8
$ws.Cells.Item(70, 10).Value = 3033.5
$ws.Cells.Item(70, 11).Value = 5144.4
$ws.Cells.Item(70, 12).Value = 9100.5
$ws.Cells.Item(70, 13).Value = -4874.4
$ws.Cells.Item(70, 14).Value = -9640.5

$ws.Cells.Item(73, 8).Value = 2300.889
$ws.Cells.Item(73, 9).Value = 1714.8
$ws.Cells.Item(73, 10).Value = 3033.5
$ws.Cells.Item(73, 11).Value = 5144.4
$ws.Cells.Item(73, 12).Value = 9100.5
$ws.Cells.Item(73, 13).Value = -4208.4
$ws.Cells.Item(73, 14).Value = -10972.5

$ws.Cells.Item(76, 8).Value = 7103.5884
$ws.Cells.Item(76, 9).Value = 6476.6
$ws.Cells.Item(76, 11).Value = 6476.6
$ws.Cells.Item(76, 13).Value = -6161.6

$ws.Cells.Item(79, 8).Value = 7103.5884
$ws.Cells.Item(79, 9).Value = 6476.6
$ws.Cells.Item(79, 11).Value = 6476.6
$ws.Cells.Item(79, 13).Value = -5384.6

$ws.Cells.Item(123, 8).Value = 58200
$ws.Cells.Item(123, 10).Value = 58200
$ws.Cells.Item(123, 12).Value = 58200
$ws.Cells.Item(123, 14).Value = -68000

$ws.Cells.Item(125, 8).Value = 4619.727
$ws.Cells.Item(125, 10).Value = 7132.6
$ws.Cells.Item(125, 12).Value = 64193.4
$ws.Cells.Item(125, 14).Value = -69113.39999999999

$ws.Cells.Item(128, 8).Value = 41950
$ws.Cells.Item(128, 10).Value = 41950
$ws.Cells.Item(128, 12).Value = 41950
$ws.Cells.Item(128, 14).Value = -51910

$ws.Cells.Item(131, 8).Value = 4267.375
$ws.Cells.Item(131, 9).Value = 2570.7273
$ws.Cells.Item(131, 11).Value = 7712.1819
$ws.Cells.Item(131, 13).Value = -2672.1819

$ws.Cells.Item(135, 8).Value = 10000674
$ws.Cells.Item(135, 9).Value = 542.35
$ws.Cells.Item(135, 11).Value = 4881.150000000001
$ws.Cells.Item(135, 13).Value = -2346.150000000001

$ws.Cells.Item(137, 8).Value = 5987.7
$ws.Cells.Item(137, 9).Value = 1973.4286
$ws.Cells.Item(137, 10).Value = 15354.333
$ws.Cells.Item(137, 11).Value = 5920.2858
$ws.Cells.Item(137, 12).Value = 46062.999
$ws.Cells.Item(137, 13).Value = -3370.2858
$ws.Cells.Item(137, 14).Value = -51162.999

$ws.Cells.Item(141, 8).Value = 6032.9165
$ws.Cells.Item(141, 9).Value = 9367.166999999999
$ws.Cells.Item(141, 11).Value = 28101.501
$ws.Cells.Item(141, 13).Value = -22921.501

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2362.4285
$ws.Cells.Item(2, 9).Value = 1784.3914
$ws.Cells.Item(2, 11).Value = 1784.3914
$ws.Cells.Item(2, 13).Value = -1671.3914

$ws.Cells.Item(61, 8).Value = 142859340
$ws.Cells.Item(61, 9).Value = 166668400
$ws.Cells.Item(61, 11).Value = 166668400
$ws.Cells.Item(61, 13).Value = -166668188

$ws.Cells.Item(74, 8).Value = 21741346
$ws.Cells.Item(74, 9).Value = 50001332
$ws.Cells.Item(74, 10).Value = 2896.7307
$ws.Cells.Item(74, 11).Value = 50001332
$ws.Cells.Item(74, 12).Value = 2896.7307
$ws.Cells.Item(74, 13).Value = -50000458
$ws.Cells.Item(74, 14).Value = -4644.7307

$ws.Cells.Item(77, 8).Value = 21741346
$ws.Cells.Item(77, 9).Value = 50001332
$ws.Cells.Item(77, 10).Value = 2896.7307
$ws.Cells.Item(77, 11).Value = 250006660
$ws.Cells.Item(77, 12).Value = 14483.6535
$ws.Cells.Item(77, 13).Value = -250002292
$ws.Cells.Item(77, 14).Value = -23219.6535

$ws.Cells.Item(103, 8).Value = 249717.5
$ws.Cells.Item(103, 10).Value = 435435
$ws.Cells.Item(103, 12).Value = 435435
$ws.Cells.Item(103, 14).Value = -437779

$ws.Cells.Item(110, 8).Value = 3610.6956
$ws.Cells.Item(110, 9).Value = 4469.7334
$ws.Cells.Item(110, 11).Value = 4469.7334
$ws.Cells.Item(110, 13).Value = -2424.7334

$ws.Cells.Item(116, 8).Value = 2362.4285
$ws.Cells.Item(116, 9).Value = 1784.3914
$ws.Cells.Item(116, 11).Value = 1784.3914
$ws.Cells.Item(116, 13).Value = 509.6086

$ws.Cells.Item(136, 8).Value = 142859340
$ws.Cells.Item(136, 9).Value = 166668400
$ws.Cells.Item(136, 11).Value = 500005200
$ws.Cells.Item(136, 13).Value = -500002650

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2362.4285
$ws.Cells.Item(3, 9).Value = 1784.3914
$ws.Cells.Item(3, 11).Value = 1784.3914
$ws.Cells.Item(3, 13).Value = -1670.3914

$ws.Cells.Item(94, 8).Value = 2286.4443
$ws.Cells.Item(94, 9).Value = 1089.1538
$ws.Cells.Item(94, 10).Value = 5399.4
$ws.Cells.Item(94, 11).Value = 1089.1538
$ws.Cells.Item(94, 12).Value = 5399.4
$ws.Cells.Item(94, 13).Value = -638.1538
$ws.Cells.Item(94, 14).Value = -6301.4

$ws.Cells.Item(107, 8).Value = 2388.2917
$ws.Cells.Item(107, 9).Value = 1485.7222
$ws.Cells.Item(107, 10).Value = 5096
$ws.Cells.Item(107, 11).Value = 1485.7222
$ws.Cells.Item(107, 12).Value = 5096
$ws.Cells.Item(107, 13).Value = 434.2778000000001
$ws.Cells.Item(107, 14).Value = -8936

$ws.Cells.Item(134, 8).Value = 2619.7097
$ws.Cells.Item(134, 9).Value = 2540.3667
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 7621.1001
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = -5086.1001
$ws.Cells.Item(134, 14).Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 8005.706
$ws.Cells.Item(86, 9).Value = 7453.5454
$ws.Cells.Item(86, 11).Value = 7453.5454
$ws.Cells.Item(86, 13).Value = -6330.5454

$ws.Cells.Item(89, 8).Value = 8005.706
$ws.Cells.Item(89, 9).Value = 7453.5454
$ws.Cells.Item(89, 11).Value = 37267.727
$ws.Cells.Item(89, 13).Value = -31651.727

$ws.Cells.Item(99, 8).Value = 3951.923
$ws.Cells.Item(99, 9).Value = 3687.7
$ws.Cells.Item(99, 10).Value = 4832.6665
$ws.Cells.Item(99, 11).Value = 3687.7
$ws.Cells.Item(99, 12).Value = 4832.6665
$ws.Cells.Item(99, 13).Value = -2189.7
$ws.Cells.Item(99, 14).Value = -7828.6665

$ws.Cells.Item(122, 8).Value = 2508192.8
$ws.Cells.Item(122, 9).Value = 1551.25
$ws.Cells.Item(122, 10).Value = 5850381.5
$ws.Cells.Item(122, 11).Value = 4653.75
$ws.Cells.Item(122, 12).Value = 17551144.5
$ws.Cells.Item(122, 13).Value = -2203.75
$ws.Cells.Item(122, 14).Value = -17556044.5

$ws.Cells.Item(126, 8).Value = 3951.923
$ws.Cells.Item(126, 9).Value = 3687.7
$ws.Cells.Item(126, 10).Value = 4832.6665
$ws.Cells.Item(126, 11).Value = 11063.1
$ws.Cells.Item(126, 12).Value = 14497.9995
$ws.Cells.Item(126, 13).Value = -8593.099999999999
$ws.Cells.Item(126, 14).Value = -19437.9995

$ws.Cells.Item(131, 8).Value = 21097.75
$ws.Cells.Item(131, 10).Value = 21097.75
$ws.Cells.Item(131, 12).Value = 21097.75
$ws.Cells.Item(131, 14).Value = -31177.75

$ws.Cells.Item(141, 8).Value = 85049.836
$ws.Cells.Item(141, 10).Value = 85049.836
$ws.Cells.Item(141, 12).Value = 85049.836
$ws.Cells.Item(141, 14).Value = -95409.836

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 9).Value = 28.666666
$ws.Cells.Item(2, 10).Value = 168.75
$ws.Cells.Item(2, 11).Value = 171.999996
$ws.Cells.Item(2, 12).Value = 1012.5
$ws.Cells.Item(2, 13).Value = -58.99999600000001
$ws.Cells.Item(2, 14).Value = -1238.5

$ws.Cells.Item(128, 8).Value = 187326
$ws.Cells.Item(128, 9).Value = 187326
$ws.Cells.Item(128, 11).Value = 561978
$ws.Cells.Item(128, 13).Value = -556998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 126126
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 13).ClearContents()

$ws.Cells.Item(66, 8).Value = 126126
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 13).ClearContents()

$ws.Cells.Item(80, 8).Value = 6025.0713
$ws.Cells.Item(80, 9).Value = 6682.875
$ws.Cells.Item(80, 11).Value = 6682.875
$ws.Cells.Item(80, 13).Value = -5684.875

$ws.Cells.Item(83, 8).Value = 6025.0713
$ws.Cells.Item(83, 9).Value = 6682.875
$ws.Cells.Item(83, 11).Value = 33414.375
$ws.Cells.Item(83, 13).Value = -28422.375

$ws.Cells.Item(132, 8).Value = 4044.2144
$ws.Cells.Item(132, 9).Value = 4044.2144
$ws.Cells.Item(132, 11).Value = 12132.6432
$ws.Cells.Item(132, 13).Value = -9602.643199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2532.3572
$ws.Cells.Item(22, 9).Value = 2281.375
$ws.Cells.Item(22, 11).Value = 2281.375
$ws.Cells.Item(22, 13).Value = -1986.375

$ws.Cells.Item(27, 8).Value = 2532.3572
$ws.Cells.Item(27, 9).Value = 2281.375
$ws.Cells.Item(27, 11).Value = 2281.375
$ws.Cells.Item(27, 13).Value = -2174.375

$ws.Cells.Item(95, 8).Value = 43606
$ws.Cells.Item(95, 10).Value = 43606
$ws.Cells.Item(95, 12).Value = 43606
$ws.Cells.Item(95, 14).Value = -49098

$ws.Cells.Item(132, 8).Value = 133343960
$ws.Cells.Item(132, 9).Value = 3228.8572
$ws.Cells.Item(132, 11).Value = 9686.571599999999
$ws.Cells.Item(132, 13).Value = -7156.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7725
$ws.Cells.Item(62, 9).Value = 5802
$ws.Cells.Item(62, 11).Value = 5802
$ws.Cells.Item(62, 13).Value = -5178

$ws.Cells.Item(64, 8).Value = 46666.668

$ws.Cells.Item(65, 8).Value = 7725
$ws.Cells.Item(65, 9).Value = 5802
$ws.Cells.Item(65, 11).Value = 29010
$ws.Cells.Item(65, 13).Value = -25890

$ws.Cells.Item(67, 8).Value = 46666.668

$ws.Cells.Item(81, 8).Value = 6571.625
$ws.Cells.Item(81, 9).Value = 6505.875
$ws.Cells.Item(81, 10).Value = 6637.375
$ws.Cells.Item(81, 11).Value = 13011.75
$ws.Cells.Item(81, 12).Value = 13274.75
$ws.Cells.Item(81, 13).Value = -11950.75
$ws.Cells.Item(81, 14).Value = -15396.75

$ws.Cells.Item(84, 8).Value = 6571.625
$ws.Cells.Item(84, 9).Value = 6505.875
$ws.Cells.Item(84, 10).Value = 6637.375
$ws.Cells.Item(84, 11).Value = 65058.75
$ws.Cells.Item(84, 12).Value = 66373.75
$ws.Cells.Item(84, 13).Value = -59754.75
$ws.Cells.Item(84, 14).Value = -76981.75

$ws.Cells.Item(96, 8).Value = 5474.6665
$ws.Cells.Item(96, 9).Value = 4324.7144
$ws.Cells.Item(96, 10).Value = 9499.5
$ws.Cells.Item(96, 11).Value = 4324.7144
$ws.Cells.Item(96, 12).Value = 9499.5
$ws.Cells.Item(96, 13).Value = -2951.7144
$ws.Cells.Item(96, 14).Value = -12245.5

$ws.Cells.Item(132, 8).Value = 1382.6316
$ws.Cells.Item(132, 9).Value = 1344.4667
$ws.Cells.Item(132, 10).Value = 1525.75
$ws.Cells.Item(132, 11).Value = 4033.4001
$ws.Cells.Item(132, 12).Value = 4577.25
$ws.Cells.Item(132, 13).Value = -1503.4001
$ws.Cells.Item(132, 14).Value = -9637.25
